$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = "'66.960.80"
$ws.Range("E2").Value = "'  -3.83%  "
$ws.Range("D3").Value = "'3.666.04"
$ws.Range("E3").Value = "'  -3.42%  "
$ws.Range("E4").Value = "'  -0.13%  "
$ws.Range("D5").Value = "'592.61"
$ws.Range("E5").Value = "'  -3.36%  "
$ws.Range("D6").Value = "'164.39"
$ws.Range("E6").Value = "'  -7.26%  "
$ws.Range("D7").Value = "'3.666.37"
$ws.Range("E7").Value = "'  -3.36%  "
$ws.Range("E8").Value = "'  -0.06%  "
$ws.Range("D9").Value = "'0.522"
$ws.Range("E9").Value = "'  -1.25%  "
$ws.Range("D10").Value = "'0.158"
$ws.Range("E10").Value = "'  -5.13%  "
$ws.Range("D11").Value = "'6.10"
$ws.Range("E11").Value = "'  -5.98%  "
$ws.Range("D12").Value = "'0.459"
$ws.Range("E12").Value = "'  -5.10%  "
$ws.Range("D13").Value = "'37.29"
$ws.Range("E13").Value = "'  -6.26%  "
$ws.Range("D14").Value = "'0.0000239"
$ws.Range("E14").Value = "'  -6.35%  "
$ws.Range("D15").Value = "'4.267.03"
$ws.Range("E15").Value = "'  -3.64%  "
$ws.Range("D16").Value = "'3.650.97"
$ws.Range("E16").Value = "'  -3.92%  "
$ws.Range("D17").Value = "'66.967.47"
$ws.Range("E17").Value = "'  -3.91%  "
$ws.Range("D18").Value = "'0.114"
$ws.Range("E18").Value = "'  -4.26%  "
$ws.Range("D19").Value = "'7.11"
$ws.Range("E19").Value = "'  -5.88%  "
$ws.Range("D20").Value = "'16.97"
$ws.Range("E20").Value = "'  +2.13%  "
$ws.Range("D21").Value = "'487.73"
$ws.Range("E21").Value = "'  -3.82%  "
$ws.Range("D22").Value = "'9.01"
$ws.Range("E22").Value = "'  -6.37%  "
$ws.Range("D23").Value = "'0.712"
$ws.Range("E23").Value = "'  -3.07%  "
$ws.Range("D24").Value = "'85.14"
$ws.Range("E24").Value = "'  -1.32%  "
$ws.Range("D25").Value = "'2.28"
$ws.Range("E25").Value = "'  -8.01%  "
$ws.Range("D26").Value = "'0.0000137"
$ws.Range("E26").Value = "'  -5.21%  "
$ws.Range("D27").Value = "'12.08"
$ws.Range("E27").Value = "'  -4.33%  "
$ws.Range("D28").Value = "'0.995"
$ws.Range("E28").Value = "'  -0.45%  "
$ws.Range("D29").Value = "'9.90"
$ws.Range("E29").Value = "'  -6.28%  "
$ws.Range("D30").Value = "'2.90"
$ws.Range("E30").Value = "'  -2.69%  "
$ws.Range("D31").Value = "'2.35"
$ws.Range("E31").Value = "'  -6.87%  "
$ws.Range("D32").Value = "'7.67"
$ws.Range("E32").Value = "'  -4.54%  "
$ws.Range("D33").Value = "'31.50"
$ws.Range("E33").Value = "'  +0.31%  "
$ws.Range("D34").Value = "'3.799.74"
$ws.Range("E34").Value = "'  -3.57%  "
$ws.Range("D35").Value = "'3.599.60"
$ws.Range("E35").Value = "'  -3.47%  "
$ws.Range("D36").Value = "'0.106"
$ws.Range("E36").Value = "'  -7.01%  "
$ws.Range("D37").Value = "'0.997"
$ws.Range("E37").Value = "'  -0.32%  "
$ws.Range("E38").Value = "'  -5.42%  "
$ws.Range("D39").Value = "'5.71"
$ws.Range("E39").Value = "'  -6.55%  "
$ws.Range("E40").Value = "'  -7.77%  "
$ws.Range("D41").Value = "'0.321"
$ws.Range("E41").Value = "'  -5.20%  "
$ws.Range("D42").Value = "'434.04"
$ws.Range("E42").Value = "'  -9.73%  "
$ws.Range("D43").Value = "'48.52"
$ws.Range("E43").Value = "'  -2.48%  "
$ws.Range("D44").Value = "'1.91"
$ws.Range("E44").Value = "'  -7.16%  "
$ws.Range("D45").Value = "'2.76"
$ws.Range("E45").Value = "'  -9.22%  "
$ws.Range("D46").Value = "'8.29"
$ws.Range("E46").Value = "'  -3.30%  "
$ws.Range("E47").Value = "'  +0.01%  "
$ws.Range("D48").Value = "'142.21"
$ws.Range("E48").Value = "'  +2.14%  "
$ws.Range("D49").Value = "'39.61"
$ws.Range("E49").Value = "'  -10.27%  "
$ws.Range("D50").Value = "'2.742.79"
$ws.Range("E50").Value = "'  -6.37%  "
$ws.Range("D51").Value = "'0.0344"
$ws.Range("E51").Value = "'  -5.15%  "
